# Correctie op filter WKT: archeologischOnderzoeksgebied zonder geometrie
# (WKT) moet weggefilterd worden, niet MET geometrie. Dit corrigeert de
# kolommen "archis_gefilterd" (E) en "verschil_ldv_min_archis_gefilterd" (F)
# voor alle rijen die afhankelijk zijn van archeologischOnderzoeksgebied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Elke entry: rijnummer, nieuwe waarde voor kolom E, nieuwe waarde voor kolom F
$updates = @(
    @(2, 108969, 97461),
    @(3, 108969, 97461),
    @(11, 103195, 104691),
    @(20, 233909, 102264),
    @(21, 233909, 102264),
    @(22, 212372, 72928),
    @(23, 212372, 72928),
    @(26, 108969, 97461),
    @(33, 1082947, 213130),
    @(40, 4035911, -777401),
    @(50, 108969, 97461),
    @(51, 108969, 97461),
    @(66, 4733, 36688),
    @(68, 908831, 438054),
    @(69, 908831, 438054),
    @(75, 4035911, -777401),
    @(76, 4035911, -777401),
    @(95, 1082947, 213130),
    @(97, 744416, -239898),
    @(98, 744416, -239898),
    @(104, 233909, 102264),
    @(106, 10798477, 3653585),
    @(111, 4865493, -310906),
    @(112, 1817816, 875954),
    @(117, 1817816, 875954),
    @(126, 615498, -98425),
    @(130, 2900763, 1089084),
    @(131, 4035911, -777401),
    @(134, 10798477, 3653585),
    @(136, 703423, 3533668),
    @(143, 452641, 2805869)
)

foreach ($u in $updates) {
    $row = $u[0]
    $eVal = $u[1]
    $fVal = $u[2]
    $ws.Cells.Item($row, 5).Value = $eVal
    $ws.Cells.Item($row, 6).Value = $fVal
}
